# Automatische test-sync: 2025-08-28 18:37:50
# Appends a new "Retour status" log row to the Logs sheet and bumps the
# matching tally on the Dashboard sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$newRow = $logs.UsedRange.Rows.Count + 1

$logs.Cells.Item($newRow, 1).Value  = "Retour status"
$logs.Cells.Item($newRow, 2).Value  = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 4).Value  = "Retour / Terugbetaling"
$logs.Cells.Item($newRow, 6).Value  = "2025-08-28 18:37:11"
$logs.Cells.Item($newRow, 7).Value  = "Ja"
$logs.Cells.Item($newRow, 8).Value  = "Nee"
$logs.Cells.Item($newRow, 9).Value  = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = $dashboard.Cells.Item(2, 2).Value() + 1

# The conditional-formatting sqref groups were anchored to the old last
# row (row 9); stretch each one down to cover the freshly appended row.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "9")
    $newRange = $logs.Range("$col" + "2:" + "$col" + $newRow)
    $fcs = $oldRange.FormatConditions
    if ($fcs.Count -gt 0) {
        $fcs.Item(1).ModifyAppliesToRange($newRange)
    }
}

